# Aider: export Bonne Nouvelle du 26.03.2015.
# Fill in the previously-empty timesheet rows on the active sheet
# ("2015.10-14") with the new entries, and set the active selection to A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "PA" block (rows 6-9) ---------------------------------------------
# Row 6: Territoires, paroisses, etc.
$ws.Range("B6").Value = 42077
$ws.Range("C6").Value = "Territoires, paroisses, etc."
$ws.Range("D6").Value = 0.05555555555555555

# Row 7: Déploiement et maintenance
$ws.Range("B7").Value = 42077
$ws.Range("C7").Value = "Déploiement et maintenance"
$ws.Range("E7").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E7").Value = 0.0763888888888889

# Row 8: Étiquettes Avery
$ws.Range("B8").Value = 42081
$ws.Range("C8").Value = "Étiquettes Avery"
$ws.Range("D8").Value = 0.013888888888888888

# Row 9: Bonne Nouvelle
$ws.Range("B9").Value = 42082
$ws.Range("C9").Value = "Bonne Nouvelle"
$ws.Range("E9").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E9").Value = 0.034722222222222224

# --- "SL" block (rows 26-35), all "MAT[CH]"/"TODO" entries --------------
$ws.Range("B26").Value = 42067
$ws.Range("C26").Value = "MAT[CH]"
$ws.Range("D26").Value = 0.125

$ws.Range("B27").Value = 42067
$ws.Range("C27").Value = "TODO"
$ws.Range("D27").Value = 0.16666666666666666

$ws.Range("B28").Value = 42069
$ws.Range("C28").Value = "TODO"
$ws.Range("D28").Value = 0.08333333333333333

$ws.Range("B29").Value = 42072
$ws.Range("C29").Value = "TODO"
$ws.Range("D29").Value = 0.2916666666666667

$ws.Range("B30").Value = 42073
$ws.Range("C30").Value = "TODO"
$ws.Range("D30").Value = 0.2916666666666667

$ws.Range("B31").Value = 42075
$ws.Range("C31").Value = "TODO"
$ws.Range("D31").Value = 0.2916666666666667

$ws.Range("B32").Value = 42076
$ws.Range("C32").Value = "TODO"
$ws.Range("D32").Value = 0.08333333333333333

$ws.Range("B33").Value = 42079
$ws.Range("C33").Value = "TODO"
$ws.Range("D33").Value = 0.08333333333333333

$ws.Range("B34").Value = 42080
$ws.Range("C34").Value = "TODO"
$ws.Range("D34").Value = 0.20833333333333334

$ws.Range("B35").Value = 42081
$ws.Range("C35").Value = "TODO"
$ws.Range("D35").Value = 0.25

# --- Selection moves to A10, matching the saved view ---------------------
$ws.Range("A10").Select()
